# ----------------------------------------------------------------------------
# Specifications.xlsx revision:
#   - renames "Feuil3" -> "Courbe d'apparition ennemis" and fills it with the
#     creature-spawn-curve parameters/formulas
#   - reworks the Scoring sheet: "LIFE" becomes "POINTS DE VIE", the VERT /
#     MARRON presence values are swapped (0 / 10), the score formulas get an
#     extra pair of parentheses, and a documentation-only text version of the
#     (future) formula is added next to each
#   - fixes a typo in the Feedback sheet ("skill" -> "kill")
# ----------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$wsScoring  = $wb.Worksheets.Item("Scoring")
$wsFeedback = $wb.Worksheets.Item("Feedback")
$wsThird    = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Rename the third sheet
# ---------------------------------------------------------------------------
$wsThird.Name = "Courbe d'apparition ennemis"

# ---------------------------------------------------------------------------
# 2. Scoring sheet
# ---------------------------------------------------------------------------
$wsScoring.Range("F2").Value = "POINTS DE VIE"
$wsScoring.Columns.Item(6).ColumnWidth = 12.7

$wsScoring.Range("C3").Value = 0
$wsScoring.Range("C4").Value = 10

$wsScoring.Range("B6").Formula = "=B3*((AVERAGE(C3:C4)/G2))"
$wsScoring.Range("C6").Value = "B3*((MOYENNE(C3:C4)/G2)*B5)"

$wsScoring.Range("B7").Formula = "=B4*((AVERAGE(C3:C4)/G2))"
$wsScoring.Range("C7").Value = "B4*((MOYENNE(C3:C4)/G2)*B5)"

# ---------------------------------------------------------------------------
# 3. Feedback sheet - fix "skill" -> "kill" typo
# ---------------------------------------------------------------------------
$wsFeedback.Range("C2").Value = "Score gagner pour le kill de l'ennemi"

# ---------------------------------------------------------------------------
# 4. Populate "Courbe d'apparition ennemis"
# ---------------------------------------------------------------------------
$wsThird.Range("A2").Value = "SCORE"
$wsThird.Range("B2").Value = 1178

$wsThird.Range("A3").Value = "COEFFICIENT DE DIFFICULTE"
$wsThird.Range("B3").Value = 1

$wsThird.Range("A4").Value = "POINTS DE VIE"
$wsThird.Range("B4").Value = 3

$wsThird.Range("A5").Value = "DUREE DE LA PARTIE"
$wsThird.Range("B5").Value = 60

$wsThird.Range("A6").Value = "NOMBRES DE CREATURES DANS L'ARENE"
$wsThird.Range("B6").Value = 1

$wsThird.Range("A7").Value = "INTERVALLE INITIAL"
$wsThird.Range("B7").Value = 3

$wsThird.Range("A9").Value = "Nombre de créatures maximum dans l'arène"
$wsThird.Range("B9").Formula = "=((B2+B3)/B5)/2"

$wsThird.Range("A10").Value = "Intervalle d'apparition des créatures"
$wsThird.Range("B10").Formula = "=B7- (((B9-B6)+B3)/(10-B4))"

# NumberFormat must be applied to B9 only after B10's formula is in place,
# otherwise the engine propagates B9's format to formulas that depend on it.
$wsThird.Range("B9").NumberFormat = "0"

$wsThird.Columns.Item(1).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 5. Restore selections per sheet (Scoring / Feedback get a selection but stay
#    inactive; the third sheet keeps the active tab, matching the workbook's
#    activeTab="2").
# ---------------------------------------------------------------------------
$wsScoring.Range("C14").Select()
$wsFeedback.Range("B16").Select()
$wsThird.Activate()
$wsThird.Range("B9").Select()
